$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had a title cell in A1 ("Arrests - Noncampus") with the
# real column headers living in row 2 and the data starting on row 3. The
# cleanup removes that stray title row entirely, promotes the header row to
# row 1, and re-cases a handful of header / label strings.

# Deleting row 1 shifts row 2 (headers) up to row 1, and all data rows up by one.
$ws.Rows.Item(1).Delete()

# Re-case the header labels now sitting in row 1.
$ws.Range("A1").Value = "Survey Year"
$ws.Range("B1").Value = "UnitID"
$ws.Range("C1").Value = "Institution Name"
$ws.Range("D1").Value = "Campus ID"
$ws.Range("E1").Value = "Campus Name"
$ws.Range("F1").Value = "Institution Size"
$ws.Range("G1").Value = "Illegal Weapons Possession"
$ws.Range("H1").Value = "Drug Law Violations"
$ws.Range("I1").Value = "Liquor Law Violations"
